# Refresh Universalis market-price snapshots and recomputed leve profit
# columns (H:N) for the rows whose prices moved since the last pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 134
$ws.Range("H134").Value = 52780
$ws.Range("J134").Value = 52780
$ws.Range("L134").Value = 52780
$ws.Range("N134").Value = -62920

# Row 137
$ws.Range("H137").Value = 457458.9
$ws.Range("I137").Value = 1404.6471
$ws.Range("K137").Value = 4213.9413
$ws.Range("M137").Value = -1663.9413

# Row 140
$ws.Range("H140").Value = 85523.75
$ws.Range("I140").Value = 70555
$ws.Range("J140").Value = 86884.55
$ws.Range("K140").Value = 70555
$ws.Range("L140").Value = 86884.55
$ws.Range("M140").Value = -65375
$ws.Range("N140").Value = -97244.55

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7647.891
$ws.Range("I32").Value = 6286.9414
$ws.Range("J32").Value = 25000
$ws.Range("K32").Value = 6286.9414
$ws.Range("L32").Value = 25000
$ws.Range("M32").Value = -5999.9414
$ws.Range("N32").Value = -25574

# Row 74
$ws.Range("H74").Value = 4745.676
$ws.Range("I74").Value = 2732.88
$ws.Range("J74").Value = 8939
$ws.Range("K74").Value = 2732.88
$ws.Range("L74").Value = 8939
$ws.Range("M74").Value = -1858.88
$ws.Range("N74").Value = -10687

# Row 77
$ws.Range("H77").Value = 4745.676
$ws.Range("I77").Value = 2732.88
$ws.Range("J77").Value = 8939
$ws.Range("K77").Value = 13664.4
$ws.Range("L77").Value = 44695
$ws.Range("M77").Value = -9296.400000000001
$ws.Range("N77").Value = -53431

# Row 122
$ws.Range("H122").Value = 2655.25
$ws.Range("I122").Value = 3890.5
$ws.Range("J122").Value = 1831.75
$ws.Range("K122").Value = 11671.5
$ws.Range("L122").Value = 5495.25
$ws.Range("M122").Value = -9221.5
$ws.Range("N122").Value = -10395.25

# Row 130
$ws.Range("H130").Value = 55323.25
$ws.Range("J130").Value = 55323.25
$ws.Range("L130").Value = 55323.25
$ws.Range("N130").Value = -65363.25

# Row 132
$ws.Range("H132").Value = 5489.9272
$ws.Range("I132").Value = 3950.4546
$ws.Range("J132").Value = 7799.136
$ws.Range("K132").Value = 11851.3638
$ws.Range("L132").Value = 23397.408
$ws.Range("M132").Value = -9321.363799999999
$ws.Range("N132").Value = -28457.408

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 4866.4243
$ws.Range("I134").Value = 4078.2856
$ws.Range("J134").Value = 9280
$ws.Range("K134").Value = 12234.8568
$ws.Range("L134").Value = 27840
$ws.Range("M134").Value = -9699.856800000001
$ws.Range("N134").Value = -32910

# Row 140
$ws.Range("H140").Value = 47074.785
$ws.Range("J140").Value = 47074.785
$ws.Range("L140").Value = 47074.785
$ws.Range("N140").Value = -57434.785

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 5054196
$ws.Range("J5").Value = 18530996
$ws.Range("L5").Value = 55592988
$ws.Range("N5").Value = -55593212

# Row 22
$ws.Range("H22").Value = 2200
$ws.Range("I22").Value = 800
$ws.Range("K22").Value = 2400
$ws.Range("M22").Value = -2231

# Row 27
$ws.Range("H27").Value = 2200
$ws.Range("I27").Value = 800
$ws.Range("K27").Value = 2400
$ws.Range("M27").Value = -2298

# Row 36
$ws.Range("I36").Value = 866.6667
$ws.Range("J36").Value = 2000
$ws.Range("K36").Value = 2600.0001
$ws.Range("L36").Value = 6000
$ws.Range("M36").Value = -2431.0001
$ws.Range("N36").Value = -6338

# Row 92
$ws.Range("H92").Value = 938.1539
$ws.Range("J92").Value = 897
$ws.Range("L92").Value = 2691
$ws.Range("N92").Value = -5187

# Row 121
$ws.Range("H121").Value = 1692.1818
$ws.Range("I121").Value = 230
$ws.Range("J121").Value = 1761.8096
$ws.Range("K121").Value = 690
$ws.Range("L121").Value = 5285.4288
$ws.Range("M121").Value = 620
$ws.Range("N121").Value = -7905.4288

# Row 131
$ws.Range("H131").Value = 777.8
$ws.Range("I131").Value = 657.8823
$ws.Range("J131").Value = 1032.625
$ws.Range("K131").Value = 1973.6469
$ws.Range("L131").Value = 3097.875
$ws.Range("M131").Value = 3066.3531
$ws.Range("N131").Value = -13177.875

# Row 132
$ws.Range("H132").Value = 2262
$ws.Range("I132").Value = 3258
$ws.Range("J132").Value = 1764
$ws.Range("K132").Value = 29322
$ws.Range("L132").Value = 15876
$ws.Range("M132").Value = -26792
$ws.Range("N132").Value = -20936

# Row 135
$ws.Range("H135").Value = 5054196
$ws.Range("J135").Value = 18530996
$ws.Range("L135").Value = 166778964
$ws.Range("N135").Value = -166784034

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 4671.8667
$ws.Range("I132").Value = 8193.200000000001
$ws.Range("J132").Value = 2911.2
$ws.Range("K132").Value = 24579.6
$ws.Range("L132").Value = 8733.599999999999
$ws.Range("M132").Value = -22049.6
$ws.Range("N132").Value = -13793.6

# Row 140
$ws.Range("H140").Value = 54565.312
$ws.Range("J140").Value = 54565.312
$ws.Range("L140").Value = 54565.312
$ws.Range("N140").Value = -64925.312

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 1400
$ws.Range("I68").Value = 1480
$ws.Range("J68").Value = 1000
$ws.Range("K68").Value = 1480
$ws.Range("L68").Value = 1000
$ws.Range("M68").Value = -731
$ws.Range("N68").Value = -2498

# Row 71
$ws.Range("H71").Value = 1400
$ws.Range("I71").Value = 1480
$ws.Range("J71").Value = 1000
$ws.Range("K71").Value = 7400
$ws.Range("L71").Value = 5000
$ws.Range("M71").Value = -3656
$ws.Range("N71").Value = -12488

# Row 112
$ws.Range("H112").Value = 31846.75
$ws.Range("J112").Value = 31846.75
$ws.Range("L112").Value = 31846.75
$ws.Range("N112").Value = -34800.75

# Row 132
$ws.Range("H132").Value = 5426.1
$ws.Range("I132").Value = 5310.3
$ws.Range("J132").Value = 5541.9
$ws.Range("K132").Value = 15930.9
$ws.Range("L132").Value = 16625.7
$ws.Range("M132").Value = -13400.9
$ws.Range("N132").Value = -21685.7

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1638.4166
$ws.Range("I132").Value = 1486.5555
$ws.Range("J132").Value = 2094
$ws.Range("K132").Value = 4459.666499999999
$ws.Range("L132").Value = 6282
$ws.Range("M132").Value = -1929.666499999999
$ws.Range("N132").Value = -11342

# Row 135
$ws.Range("H135").Value = 44370
$ws.Range("J135").Value = 44370
$ws.Range("L135").Value = 44370
$ws.Range("N135").Value = -54510
